$d = $word.ActiveDocument

# Turn on revision tracking so that replacing "are" with "feel" naturally
# splits the surrounding sentence into three runs: the text before "are",
# the newly typed "feel", and the text after it - mirroring exactly what
# happens when a human retypes a single word inside a longer run in Word.
$d.TrackRevisions = $true

$found = $d.Content
$found.Find.Execute("are underwhelmed with the help they receive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$areStart = $found.Start
$areRange = $d.Range($areStart, $areStart + 3)
$areRange.Text = "feel"

# Stop tracking and bake the single edit back in as plain text/runs (no
# <w:ins>/<w:del> markup left behind), while keeping the run split that
# was produced around the replaced word.
$d.TrackRevisions = $false
ForEach ($rev in $d.Revisions) {
    $rev.Accept()
}
